$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels for columns E, F, G
$ws.Range("E1").Value = "B2"
$ws.Range("F1").Value = "A4"
$ws.Range("G1").Value = "B3"

# Data for columns E, F, G (rows 2-21)
$data = @(
    @(9, 6, 2),
    @(4, 5, 8),
    @(6, 7, 9),
    @(3, 8, 5),
    @(4, 9, 87),
    @(3, 2, 4),
    @(7, 3, 2),
    @(9, 8, 4),
    @(1, 9, 6),
    @(1, 9, 84),
    @(2, 6, 2),
    @(9, 4, 1),
    @(0, 4, 57),
    @(7, 6, 8),
    @(4, 79, 42),
    @(0, 9, 2),
    @(0, 3, 67),
    @(5, 1, 4),
    @(3, 8, 7),
    @(58, 2, 8)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 5).Value = $data[$i][0]
    $ws.Cells.Item($r, 6).Value = $data[$i][1]
    $ws.Cells.Item($r, 7).Value = $data[$i][2]
}

# Update selection to match target state
$ws.Range("G6").Select()
